$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148; this shifts the existing rows
# 148..214 down to 149..215 (preserving their data/formatting), and
# leaves a blank row 148 ready to be populated with the new record.
$ws.Rows(148).Insert()

# Populate the newly inserted row 148 with the new weekly record.
$ws.Range("A148").Value = 7
$ws.Range("B148").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C148").Value = "Ñuble"
$ws.Range("D148").Value = 45029
$ws.Range("E148").Value = 16
$ws.Range("F148").Value = 100112040
$ws.Range("G148").Value = "Cilantro"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 100
$ws.Range("K148").Value = 1500
$ws.Range("L148").Value = 1500
$ws.Range("M148").Value = 1500
$ws.Range("N148").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O148").Value = "Provincia de Diguillín"
$ws.Range("P148").Value = 1500
$ws.Range("Q148").Value = 1
$ws.Range("R148").Value = "Hortaliza"
